# Applies scheduled-runner market-data updates to the Ultima Profits workbook.
# Updates columns H-N (currentAveragePrice*, LevePrice*, LeveProfit*) for the
# affected rows across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 136.33333
$ws.Range("J12").Value = 120
$ws.Range("L12").Value = 120
$ws.Range("N12").Value = -460
$ws.Range("H21").Value = 1033
$ws.Range("I21").Value = 1033
$ws.Range("K21").Value = 1033
$ws.Range("M21").Value = -565
$ws.Range("H23").Value = 1033
$ws.Range("I23").Value = 1033
$ws.Range("K23").Value = 1033
$ws.Range("M23").Value = -799
$ws.Range("H29").Value = 302
$ws.Range("J29").Value = 300
$ws.Range("L29").Value = 900
$ws.Range("N29").Value = -1462
$ws.Range("H40").Value = 2386.3635
$ws.Range("J40").Value = 2386.3635
$ws.Range("L40").Value = 2386.3635
$ws.Range("N40").Value = -2736.3635
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").Value = $null
$ws.Range("H47").Value = 22000
$ws.Range("J47").Value = 22000
$ws.Range("L47").Value = 22000
$ws.Range("N47").Value = -23944
$ws.Range("H101").Value = 780.9
$ws.Range("I101").Value = 1077.3334
$ws.Range("J101").Value = 336.25
$ws.Range("K101").Value = 3232.0002
$ws.Range("L101").Value = 1008.75
$ws.Range("M101").Value = -1610.0002
$ws.Range("N101").Value = -4252.75
$ws.Range("H129").Value = 1378.3188
$ws.Range("J129").Value = 1629.4642
$ws.Range("L129").Value = 4888.392599999999
$ws.Range("N129").Value = -14888.3926
$ws.Range("H132").Value = 3322.8164
$ws.Range("I132").Value = 2553.7673
$ws.Range("J132").Value = 8834.333000000001
$ws.Range("K132").Value = 7661.3019
$ws.Range("L132").Value = 26502.999
$ws.Range("M132").Value = -5131.3019
$ws.Range("N132").Value = -31562.999
$ws.Range("H137").Value = 6266357.5
$ws.Range("I137").Value = 770.13635
$ws.Range("J137").Value = 20050650
$ws.Range("K137").Value = 2310.40905
$ws.Range("L137").Value = 60151950
$ws.Range("M137").Value = 239.5909499999998
$ws.Range("N137").Value = -60157050

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1530.4546
$ws.Range("I2").Value = 1374.5714
$ws.Range("J2").Value = 1803.25
$ws.Range("K2").Value = 1374.5714
$ws.Range("L2").Value = 1803.25
$ws.Range("M2").Value = -1261.5714
$ws.Range("N2").Value = -2029.25
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = $null
$ws.Range("N3").Value = $null
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").Value = $null
$ws.Range("H32").Value = 9373.440000000001
$ws.Range("I32").Value = 7316.68
$ws.Range("J32").Value = 15543.72
$ws.Range("K32").Value = 7316.68
$ws.Range("L32").Value = 15543.72
$ws.Range("M32").Value = -7029.68
$ws.Range("N32").Value = -16117.72
$ws.Range("H35").Value = 900
$ws.Range("I35").Value = 900
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 900
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -494
$ws.Range("N35").Value = $null
$ws.Range("H38").Value = 17509.5
$ws.Range("I38").Value = 17509.5
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 17509.5
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -17042.5
$ws.Range("N38").Value = $null
$ws.Range("H41").Value = 6342.1816
$ws.Range("I41").Value = 1220.5
$ws.Range("J41").Value = 20000
$ws.Range("K41").Value = 1220.5
$ws.Range("L41").Value = 20000
$ws.Range("M41").Value = -806.5
$ws.Range("N41").Value = -20828
$ws.Range("H44").Value = 29561.625
$ws.Range("J44").Value = 29561.625
$ws.Range("L44").Value = 29561.625
$ws.Range("N44").Value = -30537.625
$ws.Range("H55").Value = 33516.8
$ws.Range("J55").Value = 33516.8
$ws.Range("L55").Value = 33516.8
$ws.Range("N55").Value = -34146.8
$ws.Range("H74").Value = 13890451
$ws.Range("I74").Value = 22728002
$ws.Range("J74").Value = 2871.1428
$ws.Range("K74").Value = 22728002
$ws.Range("L74").Value = 2871.1428
$ws.Range("M74").Value = -22727128
$ws.Range("N74").Value = -4619.1428
$ws.Range("H77").Value = 13890451
$ws.Range("I77").Value = 22728002
$ws.Range("J77").Value = 2871.1428
$ws.Range("K77").Value = 113640010
$ws.Range("L77").Value = 14355.714
$ws.Range("M77").Value = -113635642
$ws.Range("N77").Value = -23091.714
$ws.Range("H116").Value = 1530.4546
$ws.Range("I116").Value = 1374.5714
$ws.Range("J116").Value = 1803.25
$ws.Range("K116").Value = 1374.5714
$ws.Range("L116").Value = 1803.25
$ws.Range("M116").Value = 919.4286
$ws.Range("N116").Value = -6391.25
$ws.Range("H132").Value = 7144909
$ws.Range("I132").Value = 10001684
$ws.Range("J132").Value = 2972.6
$ws.Range("K132").Value = 30005052
$ws.Range("L132").Value = 8917.799999999999
$ws.Range("M132").Value = -30002522
$ws.Range("N132").Value = -13977.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1530.4546
$ws.Range("I3").Value = 1374.5714
$ws.Range("J3").Value = 1803.25
$ws.Range("K3").Value = 1374.5714
$ws.Range("L3").Value = 1803.25
$ws.Range("M3").Value = -1260.5714
$ws.Range("N3").Value = -2031.25
$ws.Range("H20").Value = 2145.65
$ws.Range("I20").Value = 2342
$ws.Range("J20").Value = 1851.125
$ws.Range("K20").Value = 2342
$ws.Range("L20").Value = 1851.125
$ws.Range("M20").Value = -2095
$ws.Range("N20").Value = -2345.125
$ws.Range("H37").Value = 315
$ws.Range("I37").Value = 315
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 315
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -178
$ws.Range("N37").Value = $null
$ws.Range("H56").Value = 50110
$ws.Range("J56").Value = 50110
$ws.Range("L56").Value = 50110
$ws.Range("N56").Value = -51588
$ws.Range("H134").Value = 3336.8867
$ws.Range("I134").Value = 2457.4333
$ws.Range("K134").Value = 7372.2999
$ws.Range("M134").Value = -4837.2999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = $null
$ws.Range("N29").Value = $null
$ws.Range("H31").Value = 7411806
$ws.Range("I31").Value = 4445.4863
$ws.Range("J31").Value = 41670850
$ws.Range("K31").Value = 4445.4863
$ws.Range("L31").Value = 41670850
$ws.Range("M31").Value = -4150.4863
$ws.Range("N31").Value = -41671440
$ws.Range("H34").Value = 7411806
$ws.Range("I34").Value = 4445.4863
$ws.Range("J34").Value = 41670850
$ws.Range("K34").Value = 4445.4863
$ws.Range("L34").Value = 41670850
$ws.Range("M34").Value = -4243.4863
$ws.Range("N34").Value = -41671254
$ws.Range("H35").Value = 1225
$ws.Range("I35").Value = 1225
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1225
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -931
$ws.Range("N35").Value = $null

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 14729.125
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 14729.125
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 44187.375
$ws.Range("M120").Value = $null
$ws.Range("N120").Value = -53863.375
$ws.Range("H137").Value = 5452.476
$ws.Range("I137").Value = 3247.5
$ws.Range("J137").Value = 8392.444
$ws.Range("K137").Value = 9742.5
$ws.Range("L137").Value = 25177.332
$ws.Range("M137").Value = -4642.5
$ws.Range("N137").Value = -35377.33199999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 466.66666
$ws.Range("I13").Value = 300
$ws.Range("J13").Value = 800
$ws.Range("K13").Value = 300
$ws.Range("L13").Value = 800
$ws.Range("M13").Value = -161
$ws.Range("N13").Value = -1078
$ws.Range("H41").Value = 885.7143
$ws.Range("I41").Value = 240
$ws.Range("K41").Value = 240
$ws.Range("M41").Value = 115
$ws.Range("H132").Value = 4367.609
$ws.Range("I132").Value = 6303.56
$ws.Range("J132").Value = 2062.9048
$ws.Range("K132").Value = 18910.68
$ws.Range("L132").Value = 6188.714399999999
$ws.Range("M132").Value = -16380.68
$ws.Range("N132").Value = -11248.7144
$ws.Range("H140").Value = 54120
$ws.Range("J140").Value = 54120
$ws.Range("L140").Value = 54120
$ws.Range("N140").Value = -64480

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4499.614
$ws.Range("I7").Value = 4224.7915
$ws.Range("J7").Value = 4829.4
$ws.Range("K7").Value = 4224.7915
$ws.Range("L7").Value = 4829.4
$ws.Range("M7").Value = -4112.7915
$ws.Range("N7").Value = -5053.4
$ws.Range("H34").Value = 11200
$ws.Range("J34").Value = 11200
$ws.Range("L34").Value = 11200
$ws.Range("N34").Value = -11544
$ws.Range("H87").Value = 32997.5
$ws.Range("J87").Value = 32997.5
$ws.Range("L87").Value = 32997.5
$ws.Range("N87").Value = -35243.5
$ws.Range("H90").Value = 32997.5
$ws.Range("J90").Value = 32997.5
$ws.Range("L90").Value = 98992.5
$ws.Range("N90").Value = -110224.5
$ws.Range("H126").Value = 4499.614
$ws.Range("I126").Value = 4224.7915
$ws.Range("J126").Value = 4829.4
$ws.Range("K126").Value = 12674.3745
$ws.Range("L126").Value = 14488.2
$ws.Range("M126").Value = -10204.3745
$ws.Range("N126").Value = -19428.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1590.2222
$ws.Range("I126").Value = 1590.2222
$ws.Range("K126").Value = 4770.6666
$ws.Range("M126").Value = -2300.6666
